$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 2 contents: becomes the area-total header row
$ws.Range("A2").Value = "SAN DIEGO AREA TOTALS"
$ws.Range("B2").Value = "Totals"

# Column width adjustments (A and B swap widths, B loses bestFit sizing)
# (values account for the 5px/MDW=6 padding Excel adds between ColumnWidth and stored width;
#  widths are quantized to whole pixels by Excel, so these are the closest achievable values
#  to the target stored widths of 23.33203125 / 12)
$ws.Columns("A").ColumnWidth = 22.5
$ws.Columns("B").ColumnWidth = 11.166666666666666

# Update selection to match the new active cell
$ws.Range("A8").Select()
